$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DeviceInfo")
$ws.Cells.Select()
